{"js": "// Update the worksheet date heading and the 25 division problems/answers.\n// Replacements are applied positionally (by paragraph / table cell) rather\n// than by text search, because several of the original problem strings\n// (e.g. \"90\u00f77=12, 6\") occur more than once in the document but map to\n// different replacement values depending on which cell they're in.\n//\n// `paragraph.insertText(text, Word.InsertLocation.replace)` is used\n// (instead of replacing on the containing Body) because it swaps only the\n// text of the existing run, preserving the run/paragraph formatting\n// (rFonts, sz, jc, etc.) already on that paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// First paragraph in the body (outside the table) holds the date line.\nparagraphs.items[0].insertText(\"2026-01-01 Thursday\", Word.InsertLocation.replace);\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only every 4th row (0, 4, 8, 12, 16 - zero based) actually holds the 5\n// problems per row; the rows in between are blank spacer rows.\nconst rowIndexes = [0, 4, 8, 12, 16];\n\nconst newValues = [\n  [\"90\u00f73=30, 0\", \"75\u00f72=37, 1\", \"15\u00f79=1, 6\", \"76\u00f74=19, 0\", \"10\u00f73=3, 1\"],\n  [\"74\u00f78=9, 2\", \"86\u00f73=28, 2\", \"35\u00f75=7, 0\", \"10\u00f76=1, 4\", \"23\u00f75=4, 3\"],\n  [\"27\u00f72=13, 1\", \"86\u00f78=10, 6\", \"45\u00f76=7, 3\", \"66\u00f74=16, 2\", \"57\u00f75=11, 2\"],\n  [\"17\u00f79=1, 8\", \"11\u00f74=2, 3\", \"90\u00f72=45, 0\", \"73\u00f73=24, 1\", \"48\u00f73=16, 0\"],\n  [\"81\u00f79=9, 0\", \"55\u00f77=7, 6\", \"90\u00f75=18, 0\", \"84\u00f75=16, 4\", \"80\u00f78=10, 0\"],\n];\n\nfor (let i = 0; i < rowIndexes.length; i++) {\n  const rowIndex = rowIndexes[i];\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const para = cell.body.paragraphs.getFirst();\n    para.insertText(newValues[i][col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and the 25 division problems/answers.\n# Replacements are applied positionally (by paragraph / table cell) rather\n# than by text search, because several of the original problem strings\n# (e.g. \"90\u00f77=12, 6\") occur more than once in the document but map to\n# different replacement values depending on which cell they're in.\n\n$d = $word.ActiveDocument\n\n# First paragraph in the body (outside the table) holds the date line.\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-01 Thursday\"\n\n$tbl = $d.Tables.Item(1)\n\n# Only every 4th row (1, 5, 9, 13, 17) actually holds the 5 problems per\n# row; the rows in between are blank spacer rows.\n$rowIndexes = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n    @(\"90\u00f73=30, 0\", \"75\u00f72=37, 1\", \"15\u00f79=1, 6\", \"76\u00f74=19, 0\", \"10\u00f73=3, 1\"),\n    @(\"74\u00f78=9, 2\", \"86\u00f73=28, 2\", \"35\u00f75=7, 0\", \"10\u00f76=1, 4\", \"23\u00f75=4, 3\"),\n    @(\"27\u00f72=13, 1\", \"86\u00f78=10, 6\", \"45\u00f76=7, 3\", \"66\u00f74=16, 2\", \"57\u00f75=11, 2\"),\n    @(\"17\u00f79=1, 8\", \"11\u00f74=2, 3\", \"90\u00f72=45, 0\", \"73\u00f73=24, 1\", \"48\u00f73=16, 0\"),\n    @(\"81\u00f79=9, 0\", \"55\u00f77=7, 6\", \"90\u00f75=18, 0\", \"84\u00f75=16, 4\", \"80\u00f78=10, 0\")\n)\n\nfor ($i = 0; $i -lt $rowIndexes.Length; $i++) {\n    $rowIndex = $rowIndexes[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $tbl.Cell($rowIndex, $col)\n        $cell.Range.Text = $newValues[$i][$col - 1]\n    }\n}\n"}
